{"js": "// Remove the trailing \"Ver no Jupiter...\" / \"\u00a9 2020 ...\" footer block (and\n// the blank paragraph that precedes it) from the end of the document body,\n// mirroring the upstream site rebuild that dropped the scraped page chrome.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the footer paragraphs by their text content so the script does not\n// depend on a brittle fixed index.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text.trim();\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") === 0) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"\\u00A9 2020\") === 0) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx !== -1 && copyrightIdx !== -1) {\n  // The blank paragraph immediately preceding \"Ver no Jupiter...\" is part of\n  // the same removed block.\n  const blankIdx = jupiterIdx - 1;\n\n  const toDelete = [];\n  if (blankIdx >= 0 && items[blankIdx].text.trim() === \"\") {\n    toDelete.push(items[blankIdx]);\n  }\n  toDelete.push(items[jupiterIdx]);\n  toDelete.push(items[copyrightIdx]);\n\n  for (const p of toDelete) {\n    p.delete();\n  }\n\n  await context.sync();\n}\n", "ps1": "# Remove the trailing \"Ver no Jupiter...\" / \"(c) 2020 ...\" footer block\n# (and the blank paragraph that precedes it) from the end of the document\n# body, mirroring the upstream site rebuild that dropped the scraped page\n# chrome.\n\n$d = $word.ActiveDocument\n$count = $d.Paragraphs.Count\n\n$jupiterIdx = -1\n$copyrightIdx = -1\n\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($jupiterIdx -eq -1 -and $t.StartsWith(\"Ver no Jupiter\")) {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Contact: luizeleno*\") {\n        $copyrightIdx = $i\n    }\n}\n\n# The blank paragraph immediately preceding \"Ver no Jupiter...\" is part of\n# the same removed block.\n$blankIdx = -1\nif ($jupiterIdx -gt 1) {\n    $prevText = $d.Paragraphs.Item($jupiterIdx - 1).Range.Text.Trim()\n    if ($prevText -eq \"\") {\n        $blankIdx = $jupiterIdx - 1\n    }\n}\n\n# Delete highest index first so earlier indices stay valid.\n$indices = @($copyrightIdx, $jupiterIdx, $blankIdx) | Where-Object { $_ -gt 0 } | Sort-Object -Descending\n\nforeach ($idx in $indices) {\n    $d.Paragraphs.Item($idx).Range.Delete()\n}\n"}
